$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 24
$ws.Range("B2").Value = 106
$ws.Range("B3").Value = 183
$ws.Range("B4").Value = 256
